$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "27.348.24"
$cell.Style = "Normal"
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.708.32"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("E4").Value = "  -0.08%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "224.24"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.74%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.5301"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -0.69%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.06620"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  -4.76%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.07676"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -0.71%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "4.512"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -2.40%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.709.38"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -0.86%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "1.943.26"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -0.97%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.5815"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -0.89%  "
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "0.0₅8183"
$cell.Style = "Normal"
$ws.Range("E16").Value = "  -1.65%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "67.69"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -0.61%  "
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "27.348.05"
$cell.Style = "Normal"
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "215.07"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -3.04%  "
$ws.Range("E20").Value = "  -0.07%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "4.632"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -2.39%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "10.42"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -2.67%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "5.985"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -2.00%  "
$ws.Range("E24").Value = "  -0.08%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "143.77"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -3.04%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "1.691"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +0.25%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "0.1204"
$cell.Style = "Normal"
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "7.246"
$cell.Style = "Normal"
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "16.22"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -2.89%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.05377"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -3.71%  "
$ws.Range("E31").Value = "  -0.85%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.482"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -2.03%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "3.433"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -0.76%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.648"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -0.86%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "2.862"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +1.20%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.9501"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -1.65%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "2.397"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -1.98%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.5861"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -1.69%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "0.01640"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -0.54%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "5.817"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -1.94%  "
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "1.050.71"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -0.57%  "
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.8438"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -1.57%  "
$ws.Range("E43").Value = "  -0.03%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "100.97"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.58%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "1.851.18"
$cell.Style = "Normal"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.0₈118"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +3.06%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "57.85"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -2.32%  "
$ws.Range("E48").Value = "  +1.82%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "1.007"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +0.39%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "8.104"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -1.13%  "
$ws.Range("E51").Value = "  -0.77%  "
